$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Add([Type]::Missing, $ws2)
$ws3.Name = "Sheet3"

$ws3.Range("A1").Value = "Aconnect"
$ws3.Range("B1").Value = "A collaboration tool for team communication and file sharing"
$ws3.Range("C1").Value = "Bcollab"
$ws3.Range("D1").Value = "A team platform that facilitates communication and file exchange"
$ws3.Range("E1").Value = 0.8

$ws3.Range("A2").Value = "Aconnect"
$ws3.Range("B2").Value = "A collaboration tool for team communication and file sharing"
$ws3.Range("C2").Value = "BExpenseManager"
$ws3.Range("D2").Value = "An application designed to manage budgets and track expenses"
$ws3.Range("E2").Value = 0.2

$ws3.Range("A3").Value = "AFiscalTrack"
$ws3.Range("B3").Value = "A financial tracking tool for monitoring budgets and expenses"
$ws3.Range("C3").Value = "Bcollab"
$ws3.Range("D3").Value = "A team platform that facilitates communication and file exchange"
$ws3.Range("E3").Value = 0.2

$ws3.Range("A4").Value = "AFiscalTrack"
$ws3.Range("B4").Value = "A financial tracking tool for monitoring budgets and expenses"
$ws3.Range("C4").Value = "BExpenseManager"
$ws3.Range("D4").Value = "An application designed to manage budgets and track expenses"
$ws3.Range("E4").Value = 0.8

# Apply the existing wrap+vcenter style (s=1) by copy-paste-format from Sheet1!A2 to A1:D4
$ws1.Range("A2").Copy()
$ws3.Range("A1:D4").PasteSpecial(-4122)

# percent-only style (new index 2), applied to E1 and E3 first
$ws3.Range("E1").NumberFormat = "0%"
$ws3.Range("E3").NumberFormat = "0%"

# percent+wrap style (new index 3): copy A1 (wrap+vcenter) style onto E2/E4, then apply percent format
$ws3.Range("A1").Copy()
$ws3.Range("E2").PasteSpecial(-4122)
$ws3.Range("E4").PasteSpecial(-4122)
$ws3.Range("E2").NumberFormat = "0%"
$ws3.Range("E4").NumberFormat = "0%"
